$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 14707456
$ws.Range("I28").Value = 20835294
$ws.Range("J28").Value = 643.2
$ws.Range("K28").Value = 20835294
$ws.Range("L28").Value = 643.2
$ws.Range("M28").Value = -20834809
$ws.Range("N28").Value = -1613.2
# Row 64
$ws.Range("H64").Value = 3062.9033
$ws.Range("I64").Value = 3039.8
$ws.Range("J64").Value = 3159.1667
$ws.Range("K64").Value = 3039.8
$ws.Range("L64").Value = 3159.1667
$ws.Range("M64").Value = -2791.8
$ws.Range("N64").Value = -3655.1667
# Row 67
$ws.Range("H67").Value = 3062.9033
$ws.Range("I67").Value = 3039.8
$ws.Range("J67").Value = 3159.1667
$ws.Range("K67").Value = 3039.8
$ws.Range("L67").Value = 3159.1667
$ws.Range("M67").Value = -2181.8
$ws.Range("N67").Value = -4875.1667
# Row 74
$ws.Range("H74").Value = 3464.2856
$ws.Range("I74").Value = 3363.6365
$ws.Range("J74").Value = 3833.3333
$ws.Range("K74").Value = 3363.6365
$ws.Range("L74").Value = 3833.3333
$ws.Range("M74").Value = -2427.6365
$ws.Range("N74").Value = -5705.3333
# Row 77
$ws.Range("H77").Value = 3464.2856
$ws.Range("I77").Value = 3363.6365
$ws.Range("J77").Value = 3833.3333
$ws.Range("K77").Value = 16818.1825
$ws.Range("L77").Value = 19166.6665
$ws.Range("M77").Value = -12138.1825
$ws.Range("N77").Value = -28526.6665
# Row 137
$ws.Range("H137").Value = 27796816
$ws.Range("I137").Value = 6250954
$ws.Range("K137").Value = 18752862
$ws.Range("M137").Value = -18750312

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 17070.908
$ws.Range("I2").Value = 18613
$ws.Range("K2").Value = 18613
$ws.Range("M2").Value = -18500
# Row 45
$ws.Range("H45").Value = 1112411
$ws.Range("I45").Value = 3334166.2
$ws.Range("J45").Value = 1533.3334
$ws.Range("K45").Value = 3334166.2
$ws.Range("L45").Value = 1533.3334
$ws.Range("M45").Value = -3333789.2
$ws.Range("N45").Value = -2287.3334
# Row 63
$ws.Range("H63").Value = 1766.6666
$ws.Range("I63").Value = 1740
$ws.Range("J63").Value = 1900
$ws.Range("K63").Value = 1740
$ws.Range("L63").Value = 1900
$ws.Range("M63").Value = -1054
$ws.Range("N63").Value = -3272
# Row 66
$ws.Range("H66").Value = 1766.6666
$ws.Range("I66").Value = 1740
$ws.Range("J66").Value = 1900
$ws.Range("K66").Value = 8700
$ws.Range("L66").Value = 9500
$ws.Range("M66").Value = -5268
$ws.Range("N66").Value = -16364
# Row 110
$ws.Range("H110").Value = 939.4
$ws.Range("I110").Value = 674.25
$ws.Range("K110").Value = 674.25
$ws.Range("M110").Value = 1370.75
# Row 116
$ws.Range("H116").Value = 17070.908
$ws.Range("I116").Value = 18613
$ws.Range("K116").Value = 18613
$ws.Range("M116").Value = -16319

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 17070.908
$ws.Range("I3").Value = 18613
$ws.Range("K3").Value = 18613
$ws.Range("M3").Value = -18499
# Row 63
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
# Row 66
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
# Row 74
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
# Row 77
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 580.6875
$ws.Range("I16").Value = 502.66666
$ws.Range("J16").Value = 814.75
$ws.Range("K16").Value = 502.66666
$ws.Range("L16").Value = 814.75
$ws.Range("M16").Value = -215.66666
$ws.Range("N16").Value = -1388.75
# Row 31
$ws.Range("H31").Value = 2818486
$ws.Range("I31").Value = 1438319.6
$ws.Range("K31").Value = 1438319.6
$ws.Range("M31").Value = -1438024.6
# Row 34
$ws.Range("H34").Value = 2818486
$ws.Range("I34").Value = 1438319.6
$ws.Range("K34").Value = 1438319.6
$ws.Range("M34").Value = -1438117.6
# Row 62
$ws.Range("H62").Value = 2859.9
$ws.Range("I62").Value = 2431.1875
$ws.Range("J62").Value = 4574.75
$ws.Range("K62").Value = 2431.1875
$ws.Range("L62").Value = 4574.75
$ws.Range("M62").Value = -1807.1875
$ws.Range("N62").Value = -5822.75
# Row 65
$ws.Range("H65").Value = 2859.9
$ws.Range("I65").Value = 2431.1875
$ws.Range("J65").Value = 4574.75
$ws.Range("K65").Value = 12155.9375
$ws.Range("L65").Value = 22873.75
$ws.Range("M65").Value = -9035.9375
$ws.Range("N65").Value = -29113.75
# Row 113
$ws.Range("H113").Value = 580.6875
$ws.Range("I113").Value = 502.66666
$ws.Range("J113").Value = 814.75
$ws.Range("K113").Value = 502.66666
$ws.Range("L113").Value = 814.75
$ws.Range("M113").Value = 1667.33334
$ws.Range("N113").Value = -5154.75
# Row 122
$ws.Range("H122").Value = 4468.8423
$ws.Range("I122").Value = 8311.111000000001
$ws.Range("J122").Value = 1010.8
$ws.Range("K122").Value = 24933.333
$ws.Range("L122").Value = 3032.4
$ws.Range("M122").Value = -22483.333
$ws.Range("N122").Value = -7932.4
# Row 132
$ws.Range("H132").Value = 2239.3142
$ws.Range("I132").Value = 1707.0834
$ws.Range("J132").Value = 3400.5454
$ws.Range("K132").Value = 5121.2502
$ws.Range("L132").Value = 10201.6362
$ws.Range("M132").Value = -2591.2502
$ws.Range("N132").Value = -15261.6362
# Row 134
$ws.Range("H134").Value = 1147610.6
$ws.Range("I134").Value = 4229.7417
$ws.Range("K134").Value = 12689.2251
$ws.Range("M134").Value = -10154.2251

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 10417.643
$ws.Range("I80").Value = 4099.3335
$ws.Range("J80").Value = 21790.6
$ws.Range("K80").Value = 4099.3335
$ws.Range("L80").Value = 21790.6
$ws.Range("M80").Value = -3101.3335
$ws.Range("N80").Value = -23786.6
# Row 83
$ws.Range("H83").Value = 10417.643
$ws.Range("I83").Value = 4099.3335
$ws.Range("J83").Value = 21790.6
$ws.Range("K83").Value = 20496.6675
$ws.Range("L83").Value = 108953
$ws.Range("M83").Value = -15504.6675
$ws.Range("N83").Value = -118937
# Row 102
$ws.Range("H102").Value = 11528.889
$ws.Range("I102").Value = 12807.5
$ws.Range("K102").Value = 12807.5
$ws.Range("M102").Value = -11185.5
# Row 132
$ws.Range("H132").Value = 33854824
$ws.Range("I132").Value = 61905388
$ws.Range("K132").Value = 185716164
$ws.Range("M132").Value = -185713634

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 583.55554
$ws.Range("I46").Value = 605
$ws.Range("J46").Value = 572.8333
$ws.Range("K46").Value = 605
$ws.Range("L46").Value = 572.8333
$ws.Range("M46").Value = -417
$ws.Range("N46").Value = -948.8333
# Row 75
$ws.Range("H75").Value = 18000
$ws.Range("J75").Value = 18000
$ws.Range("L75").Value = 18000
$ws.Range("N75").Value = -19872
# Row 78
$ws.Range("H78").Value = 18000
$ws.Range("J78").Value = 18000
$ws.Range("L78").Value = 54000
$ws.Range("N78").Value = -63360
# Row 82
$ws.Range("H82").Value = 4218.7144
$ws.Range("I82").Value = 1426.7
$ws.Range("J82").Value = 6756.909
$ws.Range("K82").Value = 1426.7
$ws.Range("L82").Value = 6756.909
$ws.Range("M82").Value = -1065.7
$ws.Range("N82").Value = -7478.909
# Row 85
$ws.Range("H85").Value = 4218.7144
$ws.Range("I85").Value = 1426.7
$ws.Range("J85").Value = 6756.909
$ws.Range("K85").Value = 1426.7
$ws.Range("L85").Value = 6756.909
$ws.Range("M85").Value = -178.7
$ws.Range("N85").Value = -9252.909
